$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 658 (pushes the existing rows 658..699 down to 659..700,
# and Excel auto-extends the used range / dimension to D700).
$ws.Rows(658).Insert()

# Populate the newly inserted row with the new daily record (2026/01/15, Thu,
# time=23, ranking=201). The leading apostrophe forces the date-like text to be
# stored as a literal string (matching column A's existing text formatting)
# instead of being auto-converted to a date serial number; resetting the style
# back to the neighbouring cell's style keeps formatting consistent.
$ws.Range("A658").Value = "'2026/01/15"
$ws.Range("A658").Style = $ws.Range("A659").Style
$ws.Range("B658").Value = "木"
$ws.Range("C658").Value = 23
$ws.Range("D658").Value = 201
